$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new data rows (34-36) following the same pattern as existing rows
$newRows = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$r = 34
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r++
}

# Reflect the selection left behind after entering the new rows (selecting the next empty row)
$ws.Range("A37:XFD1048576").Select()
